$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 26 (DiagnosticBag entries),
# shifting the rest of the table down by two rows. This makes room for the
# new ExternalFileLocation entries.
$ws.Rows("26:27").Insert()

# New row 26: MetaDslx.CodeAnalysis.Common\Diagnostic\ExternalFileLocation.cs
#             internal sealed class ExternalFileLocation
#             public sealed class ExternalFileLocation
$ws.Range("A26").Value = "MetaDslx.CodeAnalysis.Common\Diagnostic\ExternalFileLocation.cs"
$ws.Range("B26").Value = "internal sealed class ExternalFileLocation"
$ws.Range("C26").Value = "public sealed class ExternalFileLocation"

# New row 27: MetaDslx.CodeAnalysis.Common\Diagnostic\ExternalFileLocation.cs
#             internal ExternalFileLocation
#             public ExternalFileLocation
$ws.Range("A27").Value = "MetaDslx.CodeAnalysis.Common\Diagnostic\ExternalFileLocation.cs"
$ws.Range("C27").Value = "public ExternalFileLocation"
$ws.Range("B27").Value = "internal ExternalFileLocation"

# Update the view: scroll so row 16 is at the top and select A27, matching
# where the author left off editing.
$ws.Range("A27").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
